$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("clean_condensed")

# --- Correct the "moments" (column B) values ---
$corrections = @{
    2  = 21
    4  = 31
    5  = 49
    6  = 16
    7  = 20
    8  = 15
    9  = 22
    10 = 25
    11 = 24
    12 = 27
    14 = 7
    16 = 17
    17 = 24
    18 = 15
    20 = 28
    21 = 17
    22 = 41
    23 = 28
    24 = 26
    25 = 19
    26 = 21
    27 = 21
}

foreach ($row in $corrections.Keys) {
    $ws.Cells.Item($row, 2).Value = $corrections[$row]
}

# --- Re-enter the moment_rate formula across the whole range so it
#     collapses back into a single shared formula ---
$ws.Range("D2:D27").Formula = "=B2/C2"

# --- Add a "notes" column (E) with a clarifying note on row 15 ---
$ws.Range("E15").Value = "Nailogical stream; includes parts of day 7/6"
$ws.Range("E1").Value = "notes"
$ws.Columns.Item(5).ColumnWidth = 36.1

# --- Update the selection to match the saved view ---
$ws.Range("E27").Select()

$wb.Save()
